$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells where the new value looks numeric but must remain text (matches
# the source data which stores prices as literal strings, e.g. "1.00").
# Force the cell format to Text before assigning so Excel does not coerce
# the string into a number (which would drop formatting like trailing zeros).
$forceTextCells = @(
    "D5",
    "D6",
    "D8",
    "D15",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D29",
    "D30",
    "D33",
    "D36",
    "D38",
    "D39",
    "D46",
    "D47",
    "D48"
)
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values from the crypto feed refresh.
$ws.Range('D2').Value = '64.271.83'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '3.192.24'
$ws.Range('E3').Value = '  -7.12%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '563.12'
$ws.Range('E5').Value = '  -3.36%  '
$ws.Range('D6').Value = '171.05'
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '3.187.19'
$ws.Range('E9').Value = '  -7.22%  '
$ws.Range('E10').Value = '  -4.86%  '
$ws.Range('E11').Value = '  -3.73%  '
$ws.Range('E12').Value = '  -2.74%  '
$ws.Range('D13').Value = '3.739.19'
$ws.Range('E13').Value = '  -7.28%  '
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '27.63'
$ws.Range('E15').Value = '  -2.94%  '
$ws.Range('D16').Value = '64.257.00'
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').Value = '3.188.29'
$ws.Range('E18').Value = '  -7.59%  '
$ws.Range('D20').Value = '13.12'
$ws.Range('E20').Value = '  -4.35%  '
$ws.Range('D21').Value = '353.25'
$ws.Range('E21').Value = '  -3.51%  '
$ws.Range('D22').Value = '7.20'
$ws.Range('E22').Value = '  -4.43%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '69.51'
$ws.Range('E24').Value = '  -3.46%  '
$ws.Range('D25').Value = '0.506'
$ws.Range('E25').Value = '  -4.53%  '
$ws.Range('E26').Value = '  -2.12%  '
$ws.Range('D27').Value = '9.59'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').Value = '5.68'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('D33').Value = '22.17'
$ws.Range('E33').Value = '  -5.89%  '
$ws.Range('E34').Value = '  -3.92%  '
$ws.Range('E35').Value = '  -5.57%  '
$ws.Range('D36').Value = '157.16'
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('E37').Value = '  -4.96%  '
$ws.Range('D38').Value = '0.810'
$ws.Range('E38').Value = '  -7.62%  '
$ws.Range('D39').Value = '25.98'
$ws.Range('E39').Value = '  -8.56%  '
$ws.Range('E40').Value = '  -3.26%  '
$ws.Range('E41').Value = '  -3.94%  '
$ws.Range('D42').Value = '2.664.25'
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('E43').Value = '  -5.73%  '
$ws.Range('E44').Value = '  -6.32%  '
$ws.Range('E45').Value = '  -3.61%  '
$ws.Range('D46').Value = '328.47'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').Value = '38.61'
$ws.Range('E47').Value = '  -4.04%  '
$ws.Range('D48').Value = '23.66'
$ws.Range('E48').Value = '  -3.95%  '
$ws.Range('E49').Value = '  -5.72%  '
$ws.Range('E50').Value = '  -0.49%  '
